$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 80: row height 21 -> 32 (keeps customHeight) ---
$ws.Rows.Item(80).RowHeight = 32

# --- Row 147: category Hobby -> Experience ---
$ws.Range("E147").Value = "Experience"

# --- Rows 167-169: add new Drive-link URLs in G/H (column default style "2" applies automatically) ---
$ws.Range("G167").Value = "https://drive.google.com/file/d/1iRjTAlbIF1T-1Mawyl8b0wjya8UMS44Z/view?usp=sharing"
$ws.Range("H167").Value = "https://drive.google.com/file/d/1aPgghDAzmgi9_znt4qbj4Osi-bGyBrQu/view?usp=sharing"

$ws.Range("G168").Value = "https://drive.google.com/file/d/1Z8vQ9hkUCHWtQN6tkgKqPjKlqhQEplpy/view?usp=drive_link"
$ws.Range("H168").Value = "https://drive.google.com/file/d/1AbOBQUBawOiRZBLGEmfo4s9PdwlFSh0I/view?usp=sharing"

$ws.Range("G169").Value = "https://drive.google.com/file/d/1zcHVGAEUKRxnoChsPnXhMhu29HRkaKC6/view?usp=sharing"
$ws.Range("H169").Value = "https://drive.google.com/file/d/1vVQdYG6IqsLbaLPI6rfr_RNRn_5wxAcv/view?usp=sharing"

# --- Row 170: was a placeholder "Upcoming" row with a time (0.8125) -> becomes a real past
#     seminar (category Hobby) with its own Drive links; time cleared but keeps style ---
$ws.Range("E170").Value = "Hobby"
$ws.Range("G170").Value = "https://drive.google.com/file/d/1upA6nScEMGQ2b8w_5q6xKyDxnh23vB3y/view?usp=sharing"
$ws.Range("H170").Value = "https://drive.google.com/file/d/1hLfia8hRlV7Zo7WHQiS_fTPTwEOs9nrc/view?usp=sharing"
$ws.Range("I170").ClearContents()

# --- Row 171: Upcoming -> Hwealth; time cleared but keeps style ---
$ws.Range("E171").Value = "Hwealth"
$ws.Range("I171").ClearContents()

# --- Row 172: Upcoming -> Experience; time cleared but keeps style; row now gets an explicit
#     (auto) row height of 16 ---
$ws.Range("E172").Value = "Experience"
$ws.Range("I172").ClearContents()
$ws.Rows.Item(172).RowHeight = 16

# --- New row 173: the new "Upcoming" seminar entry, copying row 172's original per-column
#     formatting (date style on A, wrap style on C, time style on I) ---
$ws.Range("A172:E172").Copy()
$ws.Range("A173:E173").PasteSpecial(-4122)
$ws.Range("I172").Copy()
$ws.Range("I173").PasteSpecial(-4122)

$ws.Range("A173").Value = 45883
$ws.Range("B173").Value = "吴珊"
$ws.Range("C173").Value = "健康的未来--自我健康管理"
$ws.Range("D173").Value = "img/healthy_food.jpg"
$ws.Range("E173").Value = "Upcoming"
$ws.Range("I173").Value = 0.8125
$ws.Rows.Item(173).RowHeight = 16

# --- Selection / scroll position (best-effort; mirrors author's final cursor position) ---
$ws.Range("A166:XFD166").Select()
